# BillHubTestdata.xlsx - "Tests added in po based invoice with excelsheet"
#
# 1) POBasedInvoice sheet: refresh the PO-based-invoice test rows (2-9)
#    with a new batch of invoice numbers / quantities / base amounts.
# 2) BADashboardPage sheet: add two new "Submitting At" / "Submitting To"
#    columns (C/D), mirroring the columns already present on
#    POBasedInvoice (K/L), and widen them accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) POBasedInvoice — rows 2..9: Invoice Number (A), Base Amount (B),
#    IGST (C) and Quantity (J) all change together; columns D-I, K-L are
#    left untouched.
# ---------------------------------------------------------------------
$wsPO = $wb.Worksheets.Item("POBasedInvoice")

$poRows = @(
    @{ Row = 2; Invoice = "TESTINV97192"; Qty = "8"; Amount = "1.44" },
    @{ Row = 3; Invoice = "TESTINV04156"; Qty = "8"; Amount = "1.44" },
    @{ Row = 4; Invoice = "TESTINV17714"; Qty = "8"; Amount = "1.44" },
    @{ Row = 5; Invoice = "TESTINV73341"; Qty = "1"; Amount = "0.18" },
    @{ Row = 6; Invoice = "TESTINV34048"; Qty = "2"; Amount = "0.36" },
    @{ Row = 7; Invoice = "TESTINV25913"; Qty = "3"; Amount = "0.54" },
    @{ Row = 8; Invoice = "TESTINV06532"; Qty = "1"; Amount = "0.18" },
    @{ Row = 9; Invoice = "TESTINV09872"; Qty = "2"; Amount = "0.36" }
)

foreach ($r in $poRows) {
    $row = $r.Row

    $wsPO.Cells.Item($row, 1).Value = "'" + $r.Invoice
    $wsPO.Cells.Item($row, 1).Style = $wsPO.Cells.Item($row, 4).Style

    $wsPO.Cells.Item($row, 2).Value = "'" + $r.Qty
    $wsPO.Cells.Item($row, 2).Style = $wsPO.Cells.Item($row, 4).Style

    $wsPO.Cells.Item($row, 3).Value = "'" + $r.Amount
    $wsPO.Cells.Item($row, 3).Style = $wsPO.Cells.Item($row, 4).Style

    $wsPO.Cells.Item($row, 10).Value = "'" + $r.Qty
    $wsPO.Cells.Item($row, 10).Style = $wsPO.Cells.Item($row, 4).Style
}

$wsPO.Range("L14").Select()

# ---------------------------------------------------------------------
# 2) BADashboardPage — add "Submitting At" / "Submitting To" columns.
# ---------------------------------------------------------------------
$wsBA = $wb.Worksheets.Item("BADashboardPage")

$wsBA.Cells.Item(1, 3).Value = "Submitting At"
$wsBA.Cells.Item(1, 3).Font.Bold = $true

$wsBA.Cells.Item(1, 4).Value = "Submitting To"
$wsBA.Cells.Item(1, 4).Font.Bold = $true

$wsBA.Cells.Item(2, 3).Value = "Ahmedabad"
$wsBA.Cells.Item(2, 4).Value = "Nishant Gore"

$wsBA.Columns.Item(3).ColumnWidth = 14.665
$wsBA.Columns.Item(4).ColumnWidth = 13.665

$wsBA.Range("D5").Select()
